# Updates cryptocurrency price/volume data in the active worksheet
# Generated from the authoritative diff of D2:E51 (Price / Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.995.93"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.403.43"
$ws.Range("E3").Value = "  -3.62%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "488.05"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "154.47"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +20.10%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "2.419.76"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "6.23"
$ws.Range("E10").Value = "  +8.71%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.828.22"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "57.038.49"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "20.58"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").Value = "2.419.61"
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").Value = "325.61"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "9.95"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "58.00"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "0.404"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "2.510.70"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("E29").Value = "  -4.76%  "
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D32").Value = "149.92"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").Value = "18.52"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "3.75"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("E39").Value = "  +9.14%  "
$ws.Range("D40").Value = "34.12"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "0.994"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "0.594"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "269.25"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  -5.75%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "4.52"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").Value = "1.877.01"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").Value = "17.43"
$ws.Range("E51").Value = "  -2.57%  "
